$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Rename header K1 from "data_reperforming" to "flag_reperforming"
$ws.Range("K1").Value = "flag_reperforming"

# Add new value "N" in K2 (flag_reperforming value)
$ws.Range("K2").Value = "N"

# Update the active selection to K3, matching the saved state in the diff
$ws.Range("K3").Select()
